$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.167.55'
$ws.Range("E2").Value = '  -2.33%  '
$ws.Range("D3").Value = '1.872.19'
$ws.Range("E3").Value = '  -1.60%  '
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'307.59"
$ws.Range("E5").Value = '  -1.64%  '
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").Value = "'0.5158"
$ws.Range("E7").Value = '  +2.85%  '
$ws.Range("D8").Value = "'0.3752"
$ws.Range("E8").Value = '  -1.52%  '
$ws.Range("D9").Value = "'0.07167"
$ws.Range("E9").Value = '  -1.39%  '
$ws.Range("D10").Value = "'20.85"
$ws.Range("E10").Value = '  +0.13%  '
$ws.Range("D11").Value = "'0.8853"
$ws.Range("E11").Value = '  -2.53%  '
$ws.Range("D12").Value = '1.884.23'
$ws.Range("E12").Value = '  -0.80%  '
$ws.Range("D13").Value = "'0.07582"
$ws.Range("E13").Value = '  -0.88%  '
$ws.Range("D14").Value = "'5.338"
$ws.Range("E14").Value = '  -2.58%  '
$ws.Range("D15").Value = "'89.43"
$ws.Range("E15").Value = '  -2.02%  '
$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").Value = "'0.000008556"
$ws.Range("E17").Value = '  -1.69%  '
$ws.Range("D18").Value = "'14.21"
$ws.Range("E18").Value = '  -2.11%  '
$ws.Range("D20").Value = '27.212.44'
$ws.Range("E20").Value = '  -2.28%  '
$ws.Range("D21").Value = "'5.044"
$ws.Range("E21").Value = '  -2.21%  '
$ws.Range("D22").Value = '2.125.21'
$ws.Range("E22").Value = '  -1.13%  '
$ws.Range("E23").Value = '  -1.46%  '
$ws.Range("D24").Value = "'6.484"
$ws.Range("D25").Value = "'151.73"
$ws.Range("E25").Value = '  -1.54%  '
$ws.Range("D26").Value = "'1.856"
$ws.Range("E26").Value = '  -0.40%  '
$ws.Range("D27").Value = "'2.189"
$ws.Range("E27").Value = '  -1.93%  '
$ws.Range("E28").Value = '  -1.57%  '
$ws.Range("D29").Value = "'113.20"
$ws.Range("E29").Value = '  -1.75%  '
$ws.Range("D30").Value = "'4.755"
$ws.Range("E30").Value = '  -3.07%  '
$ws.Range("D31").Value = "'4.713"
$ws.Range("E31").Value = '  +1.62%  '
$ws.Range("D32").Value = "'0.09046"
$ws.Range("E32").Value = '  +0.86%  '
$ws.Range("D33").Value = "'0.05185"
$ws.Range("D34").Value = "'3.095"
$ws.Range("D35").Value = "'0.7581"
$ws.Range("E35").Value = '  -0.86%  '
$ws.Range("D36").Value = "'1.180"
$ws.Range("E36").Value = '  -4.10%  '
$ws.Range("D37").Value = "'0.02044"
$ws.Range("E37").Value = '  -0.54%  '
$ws.Range("D38").Value = "'2.533"
$ws.Range("E38").Value = '  -0.32%  '
$ws.Range("E39").Value = '  +0.90%  '
$ws.Range("D40").Value = "'1.084"
$ws.Range("E40").Value = '  -1.11%  '
$ws.Range("D41").Value = "'0.5433"
$ws.Range("E41").Value = '  -2.19%  '
$ws.Range("D42").Value = "'6.682"
$ws.Range("E42").Value = '  -4.16%  '
$ws.Range("D43").Value = "'115.20"
$ws.Range("E43").Value = '  +3.65%  '
$ws.Range("D44").Value = "'8.576"
$ws.Range("E44").Value = '  +1.18%  '
$ws.Range("D45").Value = "'0.1489"
$ws.Range("E45").Value = '  -1.37%  '
$ws.Range("D46").Value = "'0.4703"
$ws.Range("E46").Value = '  -1.74%  '
$ws.Range("D47").Value = "'10.21"
$ws.Range("E47").Value = '  -3.56%  '
$ws.Range("D48").Value = "'1.001"
$ws.Range("E48").Value = '  -0.16%  '
$ws.Range("D49").Value = "'1.578"
$ws.Range("E49").Value = '  -3.17%  '
$ws.Range("D50").Value = "'65.15"
$ws.Range("E50").Value = '  -3.10%  '
$ws.Range("D51").Value = "'36.57"
$ws.Range("E51").Value = '  -1.28%  '
